$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual PL_* (and swapped match-data) cells across existing rows ---
$ws.Cells.Item(2, 25).Value = -1
$ws.Cells.Item(2, 26).Value = 0.8999999999999999

$ws.Cells.Item(4, 25).Value = 0.825
$ws.Cells.Item(4, 26).Value = -1

$ws.Cells.Item(6, 25).Value = 0.95
$ws.Cells.Item(6, 26).Value = -1

$ws.Cells.Item(7, 25).Value = 0.9750000000000001
$ws.Cells.Item(7, 26).Value = -1

$ws.Cells.Item(8, 25).Value = -1
$ws.Cells.Item(8, 26).Value = 0.875

$ws.Cells.Item(9, 25).Value = 0.825
$ws.Cells.Item(9, 26).Value = -1

$ws.Cells.Item(10, 25).Value = 0.9750000000000001
$ws.Cells.Item(10, 26).Value = -1

$ws.Cells.Item(11, 2).Value = 6227815
$ws.Cells.Item(11, 5).Value = "HFX Wanderers"
$ws.Cells.Item(11, 6).Value = "Cavalry FC"
$ws.Cells.Item(11, 7).Value = 3
$ws.Cells.Item(11, 8).Value = 1
$ws.Cells.Item(11, 10).Value = 2.6
$ws.Cells.Item(11, 11).Value = 3.2
$ws.Cells.Item(11, 12).Value = 2.4
$ws.Cells.Item(11, 13).Value = 3.3
$ws.Cells.Item(11, 14).Value = 3
$ws.Cells.Item(11, 15).Value = 2.15
$ws.Cells.Item(11, 16).Value = 0.25
$ws.Cells.Item(11, 17).Value = 1.925
$ws.Cells.Item(11, 18).Value = 1.875
$ws.Cells.Item(11, 19).Value = 2.25
$ws.Cells.Item(11, 20).Value = 2
$ws.Cells.Item(11, 21).Value = 1.8
$ws.Cells.Item(11, 22).Value = 2.3
$ws.Cells.Item(11, 25).Value = 0.925
$ws.Cells.Item(11, 26).Value = -1
$ws.Cells.Item(11, 27).Value = 1
$ws.Cells.Item(11, 28).Value = -1

$ws.Cells.Item(12, 2).Value = 6240280
$ws.Cells.Item(12, 5).Value = "Atletico Ottawa"
$ws.Cells.Item(12, 6).Value = "Vancouver FC"
$ws.Cells.Item(12, 7).Value = 1
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 10).Value = 1.571
$ws.Cells.Item(12, 11).Value = 3.4
$ws.Cells.Item(12, 12).Value = 5.5
$ws.Cells.Item(12, 13).Value = 1.444
$ws.Cells.Item(12, 14).Value = 3.8
$ws.Cells.Item(12, 15).Value = 6
$ws.Cells.Item(12, 16).Value = -1.25
$ws.Cells.Item(12, 17).Value = 1.95
$ws.Cells.Item(12, 18).Value = 1.85
$ws.Cells.Item(12, 19).Value = 2.75
$ws.Cells.Item(12, 20).Value = 1.975
$ws.Cells.Item(12, 21).Value = 1.825
$ws.Cells.Item(12, 22).Value = 0.444
$ws.Cells.Item(12, 25).Value = -0.5
$ws.Cells.Item(12, 26).Value = 0.425
$ws.Cells.Item(12, 27).Value = -1
$ws.Cells.Item(12, 28).Value = 0.825

$ws.Cells.Item(13, 25).Value = 0.4875
$ws.Cells.Item(13, 26).Value = -0.5

$ws.Cells.Item(14, 25).Value = 0.925
$ws.Cells.Item(14, 26).Value = -1

$ws.Cells.Item(16, 25).Value = 0
$ws.Cells.Item(16, 26).Value = 0

$ws.Cells.Item(17, 25).Value = 0.75
$ws.Cells.Item(17, 26).Value = -1

$ws.Cells.Item(18, 25).Value = 0
$ws.Cells.Item(18, 26).Value = 0

$ws.Cells.Item(19, 25).Value = 0.75
$ws.Cells.Item(19, 26).Value = -1

$ws.Cells.Item(20, 25).Value = 1.025
$ws.Cells.Item(20, 26).Value = -1

$ws.Cells.Item(22, 25).Value = -1
$ws.Cells.Item(22, 26).Value = 0.875

$ws.Cells.Item(23, 25).Value = -1
$ws.Cells.Item(23, 26).Value = 0.825

$ws.Cells.Item(24, 25).Value = 0.875
$ws.Cells.Item(24, 26).Value = -1

$ws.Cells.Item(25, 25).Value = -1
$ws.Cells.Item(25, 26).Value = 0.95

$ws.Cells.Item(26, 25).Value = 0.8500000000000001
$ws.Cells.Item(26, 26).Value = -1

$ws.Cells.Item(27, 25).Value = -1
$ws.Cells.Item(27, 26).Value = 0.8999999999999999

$ws.Cells.Item(28, 25).Value = 0.9750000000000001
$ws.Cells.Item(28, 26).Value = -1

$ws.Cells.Item(29, 25).Value = -1
$ws.Cells.Item(29, 26).Value = 0.95

$ws.Cells.Item(31, 25).Value = 0.9750000000000001
$ws.Cells.Item(31, 26).Value = -1

$ws.Cells.Item(33, 25).Value = 0.925
$ws.Cells.Item(33, 26).Value = -1

$ws.Cells.Item(34, 25).Value = -1
$ws.Cells.Item(34, 26).Value = 0.95

$ws.Cells.Item(35, 25).Value = -1
$ws.Cells.Item(35, 26).Value = 0.95

$ws.Cells.Item(36, 25).Value = 1
$ws.Cells.Item(36, 26).Value = -1

$ws.Cells.Item(37, 25).Value = -1
$ws.Cells.Item(37, 26).Value = 0.8

$ws.Cells.Item(38, 25).Value = 0
$ws.Cells.Item(38, 26).Value = 0

$ws.Cells.Item(39, 25).Value = -1
$ws.Cells.Item(39, 26).Value = 0.925

$ws.Cells.Item(40, 25).Value = 0.7250000000000001
$ws.Cells.Item(40, 26).Value = -1

$ws.Cells.Item(41, 25).Value = -1
$ws.Cells.Item(41, 26).Value = 0.95

$ws.Cells.Item(42, 25).Value = 1.025
$ws.Cells.Item(42, 26).Value = -1

$ws.Cells.Item(45, 25).Value = -1
$ws.Cells.Item(45, 26).Value = 0.825

$ws.Cells.Item(46, 25).Value = 0.9750000000000001
$ws.Cells.Item(46, 26).Value = -1

$ws.Cells.Item(47, 25).Value = 0.95
$ws.Cells.Item(47, 26).Value = -1

$ws.Cells.Item(49, 25).Value = -1
$ws.Cells.Item(49, 26).Value = 0.925

$ws.Cells.Item(50, 25).Value = -1
$ws.Cells.Item(50, 26).Value = 0.8999999999999999

$ws.Cells.Item(51, 25).Value = 0.9750000000000001
$ws.Cells.Item(51, 26).Value = -1

$ws.Cells.Item(53, 25).Value = 0.875
$ws.Cells.Item(53, 26).Value = -1

$ws.Cells.Item(54, 25).Value = 0.4125
$ws.Cells.Item(54, 26).Value = -0.5

$ws.Cells.Item(55, 25).Value = 0.95
$ws.Cells.Item(55, 26).Value = -1

$ws.Cells.Item(57, 25).Value = 0.7749999999999999
$ws.Cells.Item(57, 26).Value = -1

$ws.Cells.Item(58, 25).Value = 0.425
$ws.Cells.Item(58, 26).Value = -0.5

$ws.Cells.Item(59, 25).Value = -1
$ws.Cells.Item(59, 26).Value = 0.75

$ws.Cells.Item(60, 25).Value = -1
$ws.Cells.Item(60, 26).Value = 0.9750000000000001

$ws.Cells.Item(61, 25).Value = -1
$ws.Cells.Item(61, 26).Value = 0.8500000000000001

$ws.Cells.Item(62, 25).Value = -1
$ws.Cells.Item(62, 26).Value = 0.8999999999999999

$ws.Cells.Item(65, 25).Value = 0.8999999999999999
$ws.Cells.Item(65, 26).Value = -1

$ws.Cells.Item(66, 25).Value = -1
$ws.Cells.Item(66, 26).Value = 0.8999999999999999

$ws.Cells.Item(69, 25).Value = -1
$ws.Cells.Item(69, 26).Value = 1

$ws.Cells.Item(70, 25).Value = -1
$ws.Cells.Item(70, 26).Value = 0.8999999999999999

$ws.Cells.Item(71, 25).Value = 0.8
$ws.Cells.Item(71, 26).Value = -1

$ws.Cells.Item(72, 25).Value = -1
$ws.Cells.Item(72, 26).Value = 0.95

$ws.Cells.Item(73, 25).Value = -1
$ws.Cells.Item(73, 26).Value = 0.8999999999999999

$ws.Cells.Item(74, 25).Value = 0.925
$ws.Cells.Item(74, 26).Value = -1

$ws.Cells.Item(75, 25).Value = 1.025
$ws.Cells.Item(75, 26).Value = -1

$ws.Cells.Item(76, 25).Value = -1
$ws.Cells.Item(76, 26).Value = 0.9750000000000001

$ws.Cells.Item(77, 25).Value = 0.375
$ws.Cells.Item(77, 26).Value = -0.5

$ws.Cells.Item(78, 25).Value = 0.875
$ws.Cells.Item(78, 26).Value = -1

$ws.Cells.Item(80, 25).Value = 0.8999999999999999
$ws.Cells.Item(80, 26).Value = -1

$ws.Cells.Item(81, 25).Value = -1
$ws.Cells.Item(81, 26).Value = 1.05

$ws.Cells.Item(82, 25).Value = -1
$ws.Cells.Item(82, 26).Value = 0.9750000000000001

$ws.Cells.Item(83, 2).Value = 6227884
$ws.Cells.Item(83, 5).Value = "Cavalry FC"
$ws.Cells.Item(83, 6).Value = "Pacific FC CA"
$ws.Cells.Item(83, 7).Value = 3
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 9).Value = "H"
$ws.Cells.Item(83, 10).Value = 2.25
$ws.Cells.Item(83, 11).Value = 3.1
$ws.Cells.Item(83, 12).Value = 2.875
$ws.Cells.Item(83, 13).Value = 2.05
$ws.Cells.Item(83, 14).Value = 3.2
$ws.Cells.Item(83, 15).Value = 3.2
$ws.Cells.Item(83, 16).Value = -0.25
$ws.Cells.Item(83, 17).Value = 1.825
$ws.Cells.Item(83, 18).Value = 1.975
$ws.Cells.Item(83, 20).Value = 1.825
$ws.Cells.Item(83, 21).Value = 1.975
$ws.Cells.Item(83, 22).Value = 1.05
$ws.Cells.Item(83, 24).Value = -1
$ws.Cells.Item(83, 25).Value = 0.825
$ws.Cells.Item(83, 26).Value = -1
$ws.Cells.Item(83, 27).Value = 0.825
$ws.Cells.Item(83, 28).Value = -1

$ws.Cells.Item(84, 2).Value = 7301364
$ws.Cells.Item(84, 5).Value = "Forge FC"
$ws.Cells.Item(84, 6).Value = "Atletico Ottawa"
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = 1
$ws.Cells.Item(84, 9).Value = "A"
$ws.Cells.Item(84, 10).Value = 1.8
$ws.Cells.Item(84, 11).Value = 3.6
$ws.Cells.Item(84, 12).Value = 3.5
$ws.Cells.Item(84, 13).Value = 1.533
$ws.Cells.Item(84, 14).Value = 3.8
$ws.Cells.Item(84, 15).Value = 5
$ws.Cells.Item(84, 16).Value = -1
$ws.Cells.Item(84, 17).Value = 1.975
$ws.Cells.Item(84, 18).Value = 1.825
$ws.Cells.Item(84, 20).Value = 1.9
$ws.Cells.Item(84, 21).Value = 1.9
$ws.Cells.Item(84, 22).Value = -1
$ws.Cells.Item(84, 24).Value = 4
$ws.Cells.Item(84, 26).Value = 0.825
$ws.Cells.Item(84, 27).Value = -1
$ws.Cells.Item(84, 28).Value = 0.8999999999999999

$ws.Cells.Item(85, 25).Value = 0.425
$ws.Cells.Item(85, 26).Value = -0.5

$ws.Cells.Item(86, 25).Value = -1
$ws.Cells.Item(86, 26).Value = 0.8

$ws.Cells.Item(87, 25).Value = -1
$ws.Cells.Item(87, 26).Value = 0.925

$ws.Cells.Item(88, 25).Value = 0.8999999999999999
$ws.Cells.Item(88, 26).Value = -1

$ws.Cells.Item(89, 25).Value = 0.8500000000000001
$ws.Cells.Item(89, 26).Value = -1

$ws.Cells.Item(90, 25).Value = 0.7749999999999999
$ws.Cells.Item(90, 26).Value = -1

$ws.Cells.Item(91, 25).Value = 0.7749999999999999
$ws.Cells.Item(91, 26).Value = -1

$ws.Cells.Item(92, 25).Value = 0.8
$ws.Cells.Item(92, 26).Value = -1

$ws.Cells.Item(93, 25).Value = 0.95
$ws.Cells.Item(93, 26).Value = -1

$ws.Cells.Item(94, 25).Value = 0.8999999999999999
$ws.Cells.Item(94, 26).Value = -1

# --- Append new row 96 (match id 94) ---
$ws.Range("A95:AB95").Copy()
$ws.Range("A96:AB96").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(96, 1).Value = 94
$ws.Cells.Item(96, 2).Value = 7802875
$ws.Cells.Item(96, 3).Value = "Canada Premier League"
$ws.Cells.Item(96, 4).Value = 45403.70833333334
$ws.Cells.Item(96, 5).Value = "York United FC"
$ws.Cells.Item(96, 6).Value = "Forge FC"
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 3
$ws.Cells.Item(96, 9).Value = "A"
$ws.Cells.Item(96, 10).Value = 2.8
$ws.Cells.Item(96, 11).Value = 3.3
$ws.Cells.Item(96, 12).Value = 2.2
$ws.Cells.Item(96, 13).Value = 3.4
$ws.Cells.Item(96, 14).Value = 3.5
$ws.Cells.Item(96, 15).Value = 1.85
$ws.Cells.Item(96, 16).Value = 0.5
$ws.Cells.Item(96, 17).Value = 1.85
$ws.Cells.Item(96, 18).Value = 1.95
$ws.Cells.Item(96, 19).Value = 2.75
$ws.Cells.Item(96, 20).Value = 1.975
$ws.Cells.Item(96, 21).Value = 1.825
$ws.Cells.Item(96, 22).Value = -1
$ws.Cells.Item(96, 23).Value = -1
$ws.Cells.Item(96, 24).Value = 0.8500000000000001
$ws.Cells.Item(96, 25).Value = -1
$ws.Cells.Item(96, 26).Value = 0.95
$ws.Cells.Item(96, 27).Value = 0.4875
$ws.Cells.Item(96, 28).Value = -0.5
